# Rename the sheet and populate it with the Selenium login test-data table
# (AbstractSeleniumTests / LoginTests fixture), matching the authored
# TestData.xlsx. Cell writes are intentionally ordered to reproduce the
# original shared-string insertion order (row 3 before row 2).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Rename the sheet.
$ws.Name = "LoginPageModel"

# 2. Header row (A1:E1) - format each cell as Text ("@") before writing so
#    the literal strings are stored as shared strings with the Text style,
#    matching the authored file's single extra cellXfs entry.
$ws.Range("A1").NumberFormat = "@"
$ws.Range("B1").NumberFormat = "@"
$ws.Range("C1").NumberFormat = "@"
$ws.Range("D1").NumberFormat = "@"
$ws.Range("E1").NumberFormat = "@"
$ws.Range("A1").Value = "Key"
$ws.Range("B1").Value = "Email"
$ws.Range("C1").Value = "Password"
$ws.Range("D1").Value = "RememberMe"
$ws.Range("E1").Value = "ExpectedError"

# 3. Row 3 - LoginShouldFailWithoutPassword (no Password supplied).
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "The Password field is required."

$ws.Range("A3").NumberFormat = "@"
$ws.Range("A3").Value = "LoginShouldFailWithoutPassword"

#    D3 holds the literal text "true" (not the boolean TRUE). Assigning the
#    string "true" straight to .Value gets auto-coerced to a boolean by the
#    engine regardless of number format, so instead produce it as a text
#    formula result in a scratch cell and paste the value in (this keeps the
#    formula's string result from being reinterpreted as a boolean on
#    write).
$ws.Range("D3").NumberFormat = "@"
$ws.Range("G1").Formula = "=""true"""
$ws.Range("G1").Copy()
$ws.Range("D3").PasteSpecial(-4163)

$ws.Range("B3").NumberFormat = "@"
$ws.Range("B3").Value = "admin@admin.com"

# 4. Row 2 - LoginShouldFailWithoutEmail (no Email supplied).
$ws.Range("A2").NumberFormat = "@"
$ws.Range("A2").Value = "LoginShouldFailWithoutEmail"

$ws.Range("C2").NumberFormat = "@"
$ws.Range("C2").Value = "1234"

$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "The Email field is required."

#    D2 also holds the literal text "true" - reuse the scratch cell's
#    clipboard contents from above.
$ws.Range("D2").NumberFormat = "@"
$ws.Range("G1").Copy()
$ws.Range("D2").PasteSpecial(-4163)
$ws.Range("G1").Clear()

# 5. Row heights as authored.
$ws.Rows("2").RowHeight = 15
$ws.Rows("3").RowHeight = 12.75

# 6. Best-effort column autosize for the data columns.
$ws.Columns("A:E").AutoFit()

# 7. Selection settles on the second row (as in the source file).
$ws.Range("A2:XFD2").Select()

# 8. Page orientation (portrait).
$ws.PageSetup.Orientation = 1
